# Minor improvement to one PowerPoint slide:
#   "  BG L1" -> "  BG L2"  (in the Content Placeholder on the
#   "Example 2" slide of "11 - Code Generation.pptx")
#
# The original run's text "  BG L1" is split in two: the leading
# two spaces stay in the original run, and "BG L1" is retyped as
# "BG L2" (forming a new run), exactly as PowerPoint does when a
# user selects "BG L1" and types "BG L2" over it.

$p = $ppt.ActivePresentation

$targetOld = "BG L1"
$targetNew = "BG L2"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if (-not $shape.HasTextFrame) { continue }

        $tr = $shape.TextFrame.TextRange
        $text = $tr.Text
        $pos = $text.IndexOf($targetOld)

        if ($pos -ge 0) {
            # COM TextRange.Characters is 1-indexed.
            $startIndex = $pos + 1
            $len = $targetOld.Length

            $sub = $tr.Characters($startIndex, $len)
            $sub.Text = $targetNew
        }
    }
}
